$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.409.26"
$ws.Range("E2").Value = "'  -2.34%  "
$ws.Range("D3").Value = "'2.577.55"
$ws.Range("E3").Value = "'  -2.84%  "
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("D5").Value = "'545.02"
$ws.Range("E5").Value = "'  +1.35%  "
$ws.Range("D6").Value = "'143.51"
$ws.Range("E6").Value = "'  -1.71%  "
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E8").Value = "'  +1.58%  "
$ws.Range("E9").Value = "'  +1.35%  "
$ws.Range("D10").Value = "'0.0997"
$ws.Range("E10").Value = "'  -3.56%  "
$ws.Range("D11").Value = "'0.139"
$ws.Range("E11").Value = "'  +3.59%  "
$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "'  -2.30%  "
$ws.Range("D13").Value = "'3.031.87"
$ws.Range("E13").Value = "'  -2.99%  "
$ws.Range("D14").Value = "'58.364.63"
$ws.Range("E14").Value = "'  -2.29%  "
$ws.Range("D15").Value = "'20.52"
$ws.Range("E15").Value = "'  -3.40%  "
$ws.Range("D16").Value = "'2.580.01"
$ws.Range("E16").Value = "'  -2.02%  "
$ws.Range("E17").Value = "'  -3.22%  "
$ws.Range("E18").Value = "'  +0.44%  "
$ws.Range("D19").Value = "'333.54"
$ws.Range("E19").Value = "'  -3.01%  "
$ws.Range("D20").Value = "'9.99"
$ws.Range("E20").Value = "'  -4.21%  "
$ws.Range("E21").Value = "'  -4.15%  "
$ws.Range("E22").Value = "'  +0.06%  "
$ws.Range("D23").Value = "'66.51"
$ws.Range("E23").Value = "'  -0.42%  "
$ws.Range("D24").Value = "'0.421"
$ws.Range("E24").Value = "'  +0.96%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "'  -0.16%  "
$ws.Range("D26").Value = "'0.158"
$ws.Range("E26").Value = "'  -5.07%  "
$ws.Range("D27").Value = "'7.04"
$ws.Range("E27").Value = "'  -4.00%  "
$ws.Range("D28").Value = "'0.0₃0734"
$ws.Range("E28").Value = "'  -2.80%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "'  +0.05%  "
$ws.Range("E30").Value = "'  -0.87%  "
$ws.Range("B31").Value = "'Monero"
$ws.Range("C31").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'154.47"
$ws.Range("E31").Value = "'  +2.59%  "
$ws.Range("B32").Value = "'Aptos"
$ws.Range("C32").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'5.88"
$ws.Range("E32").Value = "'  +0.32%  "
$ws.Range("D33").Value = "'18.83"
$ws.Range("E33").Value = "'  -0.86%  "
$ws.Range("D34").Value = "'3.88"
$ws.Range("E34").Value = "'  -3.80%  "
$ws.Range("D35").Value = "'0.847"
$ws.Range("E35").Value = "'  +0.18%  "
$ws.Range("E36").Value = "'  -4.93%  "
$ws.Range("B37").Value = "'Stacks"
$ws.Range("C37").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'1.42"
$ws.Range("E37").Value = "'  -3.75%  "
$ws.Range("B38").Value = "'Fetch.AI"
$ws.Range("C38").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'0.812"
$ws.Range("E38").Value = "'  -3.73%  "
$ws.Range("E39").Value = "'  -1.02%  "
$ws.Range("D40").Value = "'277.74"
$ws.Range("E40").Value = "'  -4.79%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "'  +0.14%  "
$ws.Range("D42").Value = "'0.593"
$ws.Range("E42").Value = "'  -2.45%  "
$ws.Range("D43").Value = "'10.62"
$ws.Range("E43").Value = "'  -1.17%  "
$ws.Range("D44").Value = "'0.0938"
$ws.Range("E44").Value = "'  -1.35%  "
$ws.Range("D45").Value = "'0.0526"
$ws.Range("E45").Value = "'  -2.52%  "
$ws.Range("D46").Value = "'18.47"
$ws.Range("E46").Value = "'  -5.45%  "
$ws.Range("D47").Value = "'0.0226"
$ws.Range("E47").Value = "'  -0.48%  "
$ws.Range("D48").Value = "'1.897.82"
$ws.Range("E48").Value = "'  -4.25%  "
$ws.Range("D49").Value = "'4.39"
$ws.Range("E49").Value = "'  -4.29%  "
$ws.Range("D50").Value = "'17.67"
$ws.Range("E50").Value = "'  -4.21%  "
$ws.Range("D51").Value = "'111.44"
